$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 16, pushing the existing rows 16-34 down to 18-36.
$ws.Range("A16:R17").EntireRow.Insert()

# Row 16: new weekly "Primera" price record (paquete 6 unidades, Región de Ñuble).
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44665
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112037
$ws.Range("G16").Value = "Cebollín"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = 650
$ws.Range("N16").Value = "$/paquete 6 unidades"
$ws.Range("O16").Value = "Región de Ñuble"
$ws.Range("P16").Value = 108
$ws.Range("Q16").Value = 6
$ws.Range("R16").Value = "Hortaliza"

# Row 17: new weekly "Segunda" price record (paquete 6 unidades, Región de Ñuble).
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44665
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112037
$ws.Range("G17").Value = "Cebollín"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = "$/paquete 6 unidades"
$ws.Range("O17").Value = "Región de Ñuble"
$ws.Range("P17").Value = 83
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = "Hortaliza"

# Apply the date number format (same style used by the rest of column D) to the new date cells.
$ws.Range("D16:D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
